{"js": "// The diff adds four paragraphs right after the first blank paragraph that\n// follows the \"Pytest with fixtures.\" line (i.e. immediately before the\n// second blank paragraph in that trailing run of blanks):\n//   1. \"\" (blank)\n//   2. \"\" (blank)\n//   3. \"Deployment: \"\n//   4. \"Heroku with postgress instance.\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet anchor = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Pytest with fixtures.\") !== -1) {\n    anchor = items[i + 1];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not locate the blank paragraph after 'Pytest with fixtures.'\");\n}\n\n// Insert each new paragraph immediately \"After\" the anchor, from last to\n// first, so they land in the correct reading order.\nanchor.insertParagraph(\"Heroku with postgress instance.\", \"After\");\nanchor.insertParagraph(\"Deployment: \", \"After\");\nanchor.insertParagraph(\"\", \"After\");\nanchor.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that contains \"Pytest with fixtures.\" so we can anchor\n# the new content on the blank paragraph immediately after it.\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Pytest with fixtures.*\") {\n        $anchorIndex = $i + 1\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate the blank paragraph after 'Pytest with fixtures.'\"\n}\n\n$anchor = $d.Paragraphs.Item($anchorIndex).Range\n\n# Insert the four new paragraphs in reverse order: each InsertParagraphAfter()\n# call on the same anchor lands immediately after it, so inserting\n# last-paragraph-first leaves everything in the correct final reading order:\n#   (anchor, blank)\n#   blank\n#   blank\n#   Deployment:\n#   Heroku with postgress instance.\n$anchor.InsertParagraphAfter()\n$d.Paragraphs.Item($anchorIndex + 1).Range.Text = \"Heroku with postgress instance.\"\n\n$anchor.InsertParagraphAfter()\n$d.Paragraphs.Item($anchorIndex + 1).Range.Text = \"Deployment: \"\n\n$anchor.InsertParagraphAfter()\n$anchor.InsertParagraphAfter()\n\nWrite-Output \"done\"\n"}
